$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value2 = 3
$ws.Range("G2").Value2 = 142.9073533333333
$ws.Range("H2").Value2 = 428.72206
$ws.Range("I2").Value2 = 0.5576664151504187
$ws.Range("J2").Value2 = 0.5576664151504188
$ws.Range("K2").Value2 = 3
$ws.Range("M2").Value2 = 21.08181366666667
$ws.Range("N2").Value2 = 63.245441
$ws.Range("O2").Value2 = 0.0571606014598545
$ws.Range("P2").Value2 = 0.0571606014598545
$ws.Range("Q2").Value2 = 3012.746194569829
$ws.Range("R2").Value2 = 27114.71575112846
$ws.Range("S2").Value2 = 0.03187654770395885
$ws.Range("T2").Value2 = 0.03187654770395885

$ws.Range("E3").Value2 = 3
$ws.Range("G3").Value2 = 142.9073533333333
$ws.Range("H3").Value2 = 428.72206
$ws.Range("I3").Value2 = 0.5576664151504187
$ws.Range("J3").Value2 = 0.5576664151504188
$ws.Range("K3").Value2 = 3
$ws.Range("M3").Value2 = 301.6001486666667
$ws.Range("N3").Value2 = 904.800446
$ws.Range("O3").Value2 = 0.8177496571571792
$ws.Range("P3").Value2 = 0.8177496571571792
$ws.Range("Q3").Value2 = 43100.87901089319
$ws.Range("R3").Value2 = 387907.9110980387
$ws.Range("S3").Value2 = 0.456031519797328
$ws.Range("T3").Value2 = 0.4560315197973281

$ws.Range("E4").Value2 = 3
$ws.Range("G4").Value2 = 142.9073533333333
$ws.Range("H4").Value2 = 428.72206
$ws.Range("I4").Value2 = 0.5576664151504187
$ws.Range("J4").Value2 = 0.5576664151504188
$ws.Range("K4").Value2 = 3
$ws.Range("M4").Value2 = 46.13524966666667
$ws.Range("N4").Value2 = 138.405749
$ws.Range("O4").Value2 = 0.1250897413829664
$ws.Range("P4").Value2 = 0.1250897413829664
$ws.Range("Q4").Value2 = 6593.066425235883
$ws.Range("R4").Value2 = 59337.59782712295
$ws.Range("S4").Value2 = 0.06975834764913183
$ws.Range("T4").Value2 = 0.06975834764913184

$ws.Range("E5").Value2 = 3
$ws.Range("G5").Value2 = 63.967809
$ws.Range("H5").Value2 = 191.903427
$ws.Range("I5").Value2 = 0.2496211559306514
$ws.Range("J5").Value2 = 0.2496211559306514
$ws.Range("K5").Value2 = 3
$ws.Range("M5").Value2 = 21.08181366666667
$ws.Range("N5").Value2 = 63.245441
$ws.Range("O5").Value2 = 0.0571606014598545
$ws.Range("P5").Value2 = 0.0571606014598545
$ws.Range("Q5").Value2 = 1348.557430002923
$ws.Range("R5").Value2 = 12137.01687002631
$ws.Range("S5").Value2 = 0.01426849541010016
$ws.Range("T5").Value2 = 0.01426849541010016

$ws.Range("E6").Value2 = 3
$ws.Range("G6").Value2 = 63.967809
$ws.Range("H6").Value2 = 191.903427
$ws.Range("I6").Value2 = 0.2496211559306514
$ws.Range("J6").Value2 = 0.2496211559306514
$ws.Range("K6").Value2 = 3
$ws.Range("M6").Value2 = 301.6001486666667
$ws.Range("N6").Value2 = 904.800446
$ws.Range("O6").Value2 = 0.8177496571571792
$ws.Range("P6").Value2 = 0.8177496571571792
$ws.Range("Q6").Value2 = 19292.70070428093
$ws.Range("R6").Value2 = 173634.3063385284
$ws.Range("S6").Value2 = 0.2041276146814689
$ws.Range("T6").Value2 = 0.2041276146814689

$ws.Range("E7").Value2 = 3
$ws.Range("G7").Value2 = 63.967809
$ws.Range("H7").Value2 = 191.903427
$ws.Range("I7").Value2 = 0.2496211559306514
$ws.Range("J7").Value2 = 0.2496211559306514
$ws.Range("K7").Value2 = 3
$ws.Range("M7").Value2 = 46.13524966666667
$ws.Range("N7").Value2 = 138.405749
$ws.Range("O7").Value2 = 0.1250897413829664
$ws.Range("P7").Value2 = 0.1250897413829664
$ws.Range("Q7").Value2 = 2951.170838844647
$ws.Range("R7").Value2 = 26560.53754960182
$ws.Range("S7").Value2 = 0.0312250458390823
$ws.Range("T7").Value2 = 0.0312250458390823

$ws.Range("E8").Value2 = 3
$ws.Range("G8").Value2 = 49.38440333333333
$ws.Range("H8").Value2 = 148.15321
$ws.Range("I8").Value2 = 0.1927124289189298
$ws.Range("J8").Value2 = 0.1927124289189298
$ws.Range("K8").Value2 = 3
$ws.Range("M8").Value2 = 21.08181366666667
$ws.Range("N8").Value2 = 63.245441
$ws.Range("O8").Value2 = 0.0571606014598545
$ws.Range("P8").Value2 = 0.0571606014598545
$ws.Range("Q8").Value2 = 1041.112789112845
$ws.Range("R8").Value2 = 9370.01510201561
$ws.Range("S8").Value2 = 0.01101555834579549
$ws.Range("T8").Value2 = 0.01101555834579549

$ws.Range("E9").Value2 = 3
$ws.Range("G9").Value2 = 49.38440333333333
$ws.Range("H9").Value2 = 148.15321
$ws.Range("I9").Value2 = 0.1927124289189298
$ws.Range("J9").Value2 = 0.1927124289189298
$ws.Range("K9").Value2 = 3
$ws.Range("M9").Value2 = 301.6001486666667
$ws.Range("N9").Value2 = 904.800446
$ws.Range("O9").Value2 = 0.8177496571571792
$ws.Range("P9").Value2 = 0.8177496571571792
$ws.Range("Q9").Value2 = 14894.34338714796
$ws.Range("R9").Value2 = 134049.0904843317
$ws.Range("S9").Value2 = 0.1575905226783821
$ws.Range("T9").Value2 = 0.1575905226783821

$ws.Range("E10").Value2 = 3
$ws.Range("G10").Value2 = 49.38440333333333
$ws.Range("H10").Value2 = 148.15321
$ws.Range("I10").Value2 = 0.1927124289189298
$ws.Range("J10").Value2 = 0.1927124289189298
$ws.Range("K10").Value2 = 3
$ws.Range("M10").Value2 = 46.13524966666667
$ws.Range("N10").Value2 = 138.405749
$ws.Range("O10").Value2 = 0.1250897413829664
$ws.Range("P10").Value2 = 0.1250897413829664
$ws.Range("Q10").Value2 = 2278.361777422699
$ws.Range("R10").Value2 = 20505.25599680429
$ws.Range("S10").Value2 = 0.02410634789475222
$ws.Range("T10").Value2 = 0.02410634789475222

